$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# Column A ("card") holds numeric-looking values but is stored as text
# throughout this sheet (e.g. "2", "10"), so force text with a leading
# apostrophe - otherwise Excel would auto-convert "10" to a number.
# Reapply the "Normal" style afterwards so the quote-prefix formatting
# doesn't leave a stray style on the cell (matches the original, which
# carries no cell-level style override on these rows).
$rows = @(3,4,5,6,7,9,10,11,12,13)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "'10"
    $cell.Style = "Normal"
}
